$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New wording for "Expected Results" (column F) - rephrased with "should"
$ws.Range("F2").Value = "Data List should contain contacts (Name, Phone No., Birth Date)"
$ws.Range("F3").Value = "Birth Date should be displayed in format YYYY-MM-DD"
$ws.Range("F4").Value = "Data List should fill the area between the top of the window and the buttons"
$ws.Range("F5").Value = "Buttons should be the same sized at the bottom of the window"
$ws.Range("F6").Value = "Objects in the window should resize when the window is resized"
$ws.Range("F7").Value = 'The "Delete" button should be disabled'
$ws.Range("F8").Value = 'The "Delete" button should becomes enabled'
$ws.Range("F9").Value = 'The selected contact should be deleted'

# "Actual Results" (column G) now filled in with what used to be the Expected Results text
$ws.Range("G2").Value = "Data List contains contacts (Name, Phone No., Birth Date)"
$ws.Range("G3").Value = "Birth Date is displayed in format YYYY-MM-DD"
$ws.Range("G4").Value = "Data List fills the area between the top of the window and the buttons"
$ws.Range("G5").Value = "Buttons are the same sized at the bottom of the window"
$ws.Range("G6").Value = "Objects in the window resize when the window is resized"
$ws.Range("G7").Value = 'The "Delete" button is disabled'
$ws.Range("G8").Value = 'The "Delete" button becomes enabled'
$ws.Range("G9").Value = "The selected contact is deleted"

# Column G7/G8 previously used a "blank" style lacking the wrap/centre alignment that the
# other Actual Results cells use; match it up with the already-used formatting (same as G2)
# without introducing a brand new style entry.
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G7").PasteSpecial(-4122) | Out-Null
$ws.Range("G8").PasteSpecial(-4122) | Out-Null

# Test Status Pass/Fail (column H) - mark every test case as Pass
$ws.Range("H2").Value = "Pass"
$ws.Range("H3").Value = "Pass"
$ws.Range("H4").Value = "Pass"
$ws.Range("H5").Value = "Pass"
$ws.Range("H6").Value = "Pass"
$ws.Range("H7").Value = "Pass"
$ws.Range("H8").Value = "Pass"
$ws.Range("H9").Value = "Pass"

$excel.CutCopyMode = 0

# Update the saved selection/scroll position to match the authored workbook
[void]$ws.Range("H4").Select()
